$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-12-01"

# Update the header label cell (I1) text that references the "through" date
$ws.Range("I1").Value = "2022 (through 12-01)"

# Add new December 2022 data point
$ws.Range("I13").Value = 6

# Update the Total row for the 2022 column to include the new value
$ws.Range("I14").Value = 1522
